$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the timestamp changes ---
$ws.Range("D2").Value = "2025-03-06 08:54:05"

# --- Row 3: timestamp changes, execution_time changes ---
$ws.Range("D3").Value = "2025-03-06 08:54:05"
$ws.Range("M3").Value = 0.003

# --- Row 4: indices shift, timestamp changes, url/api_endpoint swap to memo ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "2025-03-06 08:54:05"
$ws.Range("F4").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G4").Value = "/api/v1/memo/21"

# --- Row 5: indices shift, timestamp changes, url/api_endpoint swap to resource,
#            execution_time/user_type/data_type/data_valid/seq_valid change ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "2025-03-06 08:54:05"
$ws.Range("F5").Value = "http://49.234.6.241:5230/api/v1/resource/16"
$ws.Range("G5").Value = "/api/v1/resource/16"
$ws.Range("M5").Value = 0.004
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $true
$ws.Range("Q5").Value = $true

# --- Row 6: brand-new row ---
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = "2025-03-06 08:54:05"
$ws.Range("E6").Value = "DELETE"
$ws.Range("F6").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G6").Value = "/api/v1/memo/21"
$ws.Range("H6").Value = "{}"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 200
$ws.Range("M6").Value = 0.003
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $true
$ws.Range("Q6").Value = $false

# --- Row 7: brand-new row ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = "2025-03-06 08:54:05"
$ws.Range("E7").Value = "DELETE"
$ws.Range("F7").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G7").Value = "/api/v1/memo/21"
$ws.Range("H7").Value = "{}"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 5
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 0.004
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = $false
